$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Submit" button used the same token/id as another submit control.
# Give the age-confirmation submit button its own token ("submit_age") by
# duplicating the existing "submit" row (row 21: token_data!A21:D21) as a
# new row, inserted just above it, keeping the same label/value/description
# strings but with the new unique token name in column A.

$ws.Rows(21).Insert()

$ws.Range("A21").Value = "submit_age"
$ws.Range("B21").Value = $ws.Range("B22").Text
$ws.Range("C21").Value = $ws.Range("C22").Text
$ws.Range("D21").Value = $ws.Range("D22").Text

# New cells otherwise pick up the column's default style (wrap-text) via
# the <cols> style; match the rest of the table's unstyled cells.
$ws.Range("A21:D21").Style = "Normal"

# Keep the named range / used-range in sync with the extra row.
$names = $wb.Names
$n = $names.Item(1)
$n.RefersTo = "=token_data!`$A`$1:`$D`$86"

# Match the author's final selection in the sheet.
$ws.Range("B20").Select()
